{"js": "// Each row of the practice-sheet table holds a \"NN\u00d7NN=\" problem; update the\n// 25 populated cells to the new operands while leaving everything else\n// (fonts, sizes, empty rows, the date header, ...) untouched.\nconst replacements = [\n  [\"19\u00d777=\", \"28\u00d772=\"],\n  [\"67\u00d772=\", \"95\u00d719=\"],\n  [\"94\u00d724=\", \"64\u00d719=\"],\n  [\"95\u00d784=\", \"42\u00d776=\"],\n  [\"83\u00d794=\", \"38\u00d711=\"],\n  [\"13\u00d783=\", \"39\u00d754=\"],\n  [\"92\u00d795=\", \"77\u00d781=\"],\n  [\"12\u00d792=\", \"72\u00d782=\"],\n  [\"22\u00d787=\", \"17\u00d749=\"],\n  [\"25\u00d759=\", \"66\u00d761=\"],\n  [\"50\u00d767=\", \"34\u00d715=\"],\n  [\"44\u00d799=\", \"34\u00d754=\"],\n  [\"64\u00d726=\", \"43\u00d771=\"],\n  [\"81\u00d750=\", \"60\u00d795=\"],\n  [\"27\u00d780=\", \"75\u00d743=\"],\n  [\"13\u00d778=\", \"81\u00d743=\"],\n  [\"46\u00d759=\", \"55\u00d764=\"],\n  [\"55\u00d720=\", \"79\u00d788=\"],\n  [\"52\u00d730=\", \"55\u00d771=\"],\n  [\"13\u00d765=\", \"54\u00d784=\"],\n  [\"70\u00d745=\", \"26\u00d788=\"],\n  [\"56\u00d749=\", \"65\u00d722=\"],\n  [\"20\u00d768=\", \"21\u00d754=\"],\n  [\"74\u00d715=\", \"32\u00d724=\"],\n  [\"34\u00d766=\", \"66\u00d712=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const found = context.document.body.search(oldText, { matchCase: true });\n  found.load(\"text\");\n  await context.sync();\n\n  for (let i = 0; i < found.items.length; i++) {\n    found.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Each row of the practice-sheet table holds a \"NN\u00d7NN=\" problem; update the\n# 25 populated cells to the new operands while leaving everything else\n# (fonts, sizes, empty rows, the date header, ...) untouched.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"19\u00d777=\"; New = \"28\u00d772=\" },\n    @{ Old = \"67\u00d772=\"; New = \"95\u00d719=\" },\n    @{ Old = \"94\u00d724=\"; New = \"64\u00d719=\" },\n    @{ Old = \"95\u00d784=\"; New = \"42\u00d776=\" },\n    @{ Old = \"83\u00d794=\"; New = \"38\u00d711=\" },\n    @{ Old = \"13\u00d783=\"; New = \"39\u00d754=\" },\n    @{ Old = \"92\u00d795=\"; New = \"77\u00d781=\" },\n    @{ Old = \"12\u00d792=\"; New = \"72\u00d782=\" },\n    @{ Old = \"22\u00d787=\"; New = \"17\u00d749=\" },\n    @{ Old = \"25\u00d759=\"; New = \"66\u00d761=\" },\n    @{ Old = \"50\u00d767=\"; New = \"34\u00d715=\" },\n    @{ Old = \"44\u00d799=\"; New = \"34\u00d754=\" },\n    @{ Old = \"64\u00d726=\"; New = \"43\u00d771=\" },\n    @{ Old = \"81\u00d750=\"; New = \"60\u00d795=\" },\n    @{ Old = \"27\u00d780=\"; New = \"75\u00d743=\" },\n    @{ Old = \"13\u00d778=\"; New = \"81\u00d743=\" },\n    @{ Old = \"46\u00d759=\"; New = \"55\u00d764=\" },\n    @{ Old = \"55\u00d720=\"; New = \"79\u00d788=\" },\n    @{ Old = \"52\u00d730=\"; New = \"55\u00d771=\" },\n    @{ Old = \"13\u00d765=\"; New = \"54\u00d784=\" },\n    @{ Old = \"70\u00d745=\"; New = \"26\u00d788=\" },\n    @{ Old = \"56\u00d749=\"; New = \"65\u00d722=\" },\n    @{ Old = \"20\u00d768=\"; New = \"21\u00d754=\" },\n    @{ Old = \"74\u00d715=\"; New = \"32\u00d724=\" },\n    @{ Old = \"34\u00d766=\"; New = \"66\u00d712=\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.Execute($r.Old, $false, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)\n}\n"}
